$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '62.378.12'
$ws.Range("E2").Value = '  -0.48%  '
Set-TextValue $ws.Range("D3") '2.452.50'
$ws.Range("E3").Value = '  +0.32%  '
Set-TextValue $ws.Range("D4") '0.998'
Set-TextValue $ws.Range("D5") '576.11'
$ws.Range("E5").Value = '  +1.05%  '
Set-TextValue $ws.Range("D6") '143.97'
$ws.Range("E6").Value = '  -0.55%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +0.03%  '
Set-TextValue $ws.Range("D9") '2.447.26'
$ws.Range("E9").Value = '  +0.28%  '
$ws.Range("E10").Value = '  -0.47%  '
$ws.Range("E11").Value = '  +2.53%  '
Set-TextValue $ws.Range("D12") '5.21'
$ws.Range("E12").Value = '  -0.53%  '
$ws.Range("E13").Value = '  -3.08%  '
Set-TextValue $ws.Range("D14") '26.32'
$ws.Range("E14").Value = '  -2.91%  '
Set-TextValue $ws.Range("D15") '0.0000176'
$ws.Range("E15").Value = '  -0.02%  '
Set-TextValue $ws.Range("D16") '2.896.28'
$ws.Range("E16").Value = '  +0.44%  '
Set-TextValue $ws.Range("D17") '62.052.92'
$ws.Range("E17").Value = '  -0.86%  '
Set-TextValue $ws.Range("D18") '2.445.94'
$ws.Range("E18").Value = '  +1.29%  '
Set-TextValue $ws.Range("D19") '10.88'
$ws.Range("E19").Value = '  -3.08%  '
Set-TextValue $ws.Range("D20") '7.13'
$ws.Range("E20").Value = '  -1.59%  '
Set-TextValue $ws.Range("D21") '328.66'
$ws.Range("E21").Value = '  +0.38%  '
$ws.Range("E22").Value = '  -1.39%  '
$ws.Range("E23").Value = '  -7.12%  '
$ws.Range("E24").Value = '  +0.06%  '
Set-TextValue $ws.Range("D25") '65.67'
$ws.Range("E25").Value = '  +0.51%  '
Set-TextValue $ws.Range("D26") '9.23'
$ws.Range("E26").Value = '  +2.38%  '
Set-TextValue $ws.Range("D27") '591.70'
$ws.Range("E27").Value = '  -5.48%  '
Set-TextValue $ws.Range("D28") '2.575.25'
$ws.Range("E28").Value = '  +0.36%  '
Set-TextValue $ws.Range("D29") '0.0₃0957'
$ws.Range("E29").Value = '  -4.23%  '
$ws.Range("E30").Value = '  -0.14%  '
$ws.Range("E31").Value = '  -4.33%  '
Set-TextValue $ws.Range("D32") '8.03'
$ws.Range("E32").Value = '  -1.54%  '
$ws.Range("E33").Value = '  +0.79%  '
$ws.Range("E34").Value = '  -0.89%  '
Set-TextValue $ws.Range("D35") '4.94'
$ws.Range("E35").Value = '  -3.68%  '
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("E37").Value = '  -3.44%  '
$ws.Range("E38").Value = '  +0.22%  '
Set-TextValue $ws.Range("D39") '151.87'
$ws.Range("E39").Value = '  +3.61%  '
Set-TextValue $ws.Range("D40") '5.35'
$ws.Range("E40").Value = '  +0.44%  '
Set-TextValue $ws.Range("D41") '18.40'
$ws.Range("E41").Value = '  -2.15%  '
$ws.Range("E42").Value = '  -2.28%  '
Set-TextValue $ws.Range("D43") '42.63'
$ws.Range("E43").Value = '  +0.93%  '
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("E45").Value = '  -4.39%  '
Set-TextValue $ws.Range("D46") '142.05'
$ws.Range("E46").Value = '  -3.22%  '
Set-TextValue $ws.Range("D47") '3.64'
$ws.Range("E47").Value = '  -3.21%  '
$ws.Range("E48").Value = '  +1.12%  '
Set-TextValue $ws.Range("D49") '0.0522'
$ws.Range("E49").Value = '  -1.55%  '
Set-TextValue $ws.Range("D50") '0.0₆0245'
$ws.Range("E50").Value = '  +7.84%  '
Set-TextValue $ws.Range("D51") '19.78'
$ws.Range("E51").Value = '  -4.52%  '
